# Update "想去人数" (want-to-go count) figures across sheets, as produced by
# the latest gh-pages data regeneration (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2826
$ws1.Range("F12").Value = 328
$ws1.Range("F21").Value = 544
$ws1.Range("F24").Value = 108

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 1007
$ws2.Range("F24").Value = 4056
$ws2.Range("F28").Value = 145

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value = 1495

# 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1495
$ws4.Range("F12").Value = 2826
$ws4.Range("F18").Value = 328
$ws4.Range("F27").Value = 544
$ws4.Range("F39").Value = 145
